# Worked on reversing words, in a string
# Adds a new day entry (row 13) to the workload tracking sheet, mirroring
# the style/format of the preceding rows, and records a new task
# ("Reverse words") together with its time spent.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 is a new day, one day after the last existing entry (row 12).
$lastRow = 12
$newRow = 13

# Duplicate the formatting of the previous row onto the new row first,
# so the new cells inherit the same styles (date format, borders, etc.)
$ws.Range("A$lastRow`:E$lastRow").Copy()
$ws.Range("A$newRow`:E$newRow").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Fill in the new row's values.
$ws.Range("A$newRow").Value = 41320                # 2013-02-15
$ws.Range("B$newRow").Value = "0.5h"
$ws.Range("C$newRow").Value = "0.75H"
$ws.Range("D$newRow").Value = "0h"
$ws.Range("E$newRow").Value = "Reverse words"

# Reflect the edit in the current selection.
$ws.Range("E$newRow").Select()

$wb.Save()
